$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: a follow-up calendar entry appended below the existing two rows.
$ws.Range("A4").Value = "TEST/EVENT"
$ws.Range("B4").Value = "MYMY"
$ws.Range("C4").Value = "Hello"
$ws.Range("D4").Value = "Tutorial"

$ws.Range("E4").Value = 45884
$ws.Range("E4").NumberFormat = "d-mmm"

$ws.Range("F4").NumberFormat = "h:mm"

$ws.Range("G4").Value = 45885
$ws.Range("G4").NumberFormat = "d-mmm"

$ws.Range("I4").Value = "Sydney, Australia"
$ws.Range("J4").Value = "Sydney, Australia"
$ws.Range("M4").Value = "Transparent"

# Move the active selection as recorded in the saved view state.
[void]$ws.Range("I5").Select()
